$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data
# Row 2
$ws.Range("D2").Value = "'35.007.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "'1.814.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.51%  "

# Row 5
$ws.Range("D5").Value = "'233.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.19%  "

# Row 6
$ws.Range("D6").Value = "'0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "

# Row 7
$ws.Range("E7").Value = "  +0.58%  "

# Row 8
$ws.Range("D8").Value = "'40.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.06%  "

# Row 9
$ws.Range("D9").Value = "'0.318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.54%  "

# Row 10
$ws.Range("E10").Value = "  +1.31%  "

# Row 11
$ws.Range("D11").Value = "'0.1000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "

# Row 12
$ws.Range("D12").Value = "'2.064.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.810.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.36%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.34%  "

# Row 16
$ws.Range("E16").Value = "  +2.83%  "

# Row 17
$ws.Range("D17").Value = "'34.957.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.23%  "

# Row 18
$ws.Range("D18").Value = "'69.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.92%  "

# Row 19
$ws.Range("D19").Value = "'0.0₃0790"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "

# Row 20
$ws.Range("D20").Value = "'237.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "

# Row 21
$ws.Range("D21").Value = "'11.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "

# Row 22
$ws.Range("D22").Value = "'4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.42%  "

# Row 23
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
$ws.Range("E24").Value = "  +5.37%  "

# Row 25
$ws.Range("D25").Value = "'172.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("D26").Value = "'7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("D27").Value = "'17.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "

# Row 28
$ws.Range("E28").Value = "  -0.90%  "

# Row 29
$ws.Range("D29").Value = "'1.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +33.29%  "

# Row 30
$ws.Range("D30").Value = "'1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "

# Row 31
$ws.Range("D31").Value = "'3.339.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +37.43%  "

# Row 32
$ws.Range("E32").Value = "  +6.62%  "

# Row 33
$ws.Range("D33").Value = "'3.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "

# Row 34
$ws.Range("E34").Value = "  +1.50%  "

# Row 35
$ws.Range("E35").Value = "  -3.09%  "

# Row 36
$ws.Range("D36").Value = "'93.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.87%  "

# Row 37
$ws.Range("E37").Value = "  +7.23%  "

# Row 38
$ws.Range("E38").Value = "  +3.22%  "

# Row 39
$ws.Range("E39").Value = "  +1.25%  "

# Row 40
$ws.Range("E40").Value = "  +5.48%  "

# Row 41
$ws.Range("D41").Value = "'1.306.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

# Row 42
$ws.Range("D42").Value = "'0.988"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.06%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'14.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.47%  "

# Row 45
$ws.Range("E45").Value = "  +0.89%  "

# Row 46
$ws.Range("D46").Value = "'2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "

# Row 47
$ws.Range("D47").Value = "'6.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.63%  "

# Row 48
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("D49").Value = "'1.987.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$ws.Range("E50").Value = "  +0.53%  "

# Row 51
$ws.Range("E51").Value = "  +5.36%  "
